$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row: D, M, N, O, P, Q, R, S, T
# Derived from the weekly reshuffle of the dataset (rows re-ordered by date).
$rows = @(
    @{ Row = 2; D = 44315; M = 65; N = 14000; O = 14000; P = 14000; Q = '$/caja 14 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 14 },
    @{ Row = 3; D = 44312; M = 68; N = 14000; O = 14000; P = 14000; Q = '$/caja 14 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 14 },
    @{ Row = 4; D = 44588; M = 85; N = 19000; O = 20000; P = 19529; Q = '$/caja 14 kilos granel'; R = 'Provincia de Limarí'; S = 1395; T = 14 },
    @{ Row = 5; D = 44259; M = 80; N = 12000; O = 12000; P = 12000; Q = '$/caja 15 kilos empedrada'; R = 'Provincia de Limarí'; S = 800; T = 15 },
    @{ Row = 6; D = 44313; M = 36; N = 14000; O = 14000; P = 14000; Q = '$/caja 14 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 14 },
    @{ Row = 7; D = 44630; M = 75; N = 15000; O = 15000; P = 15000; Q = '$/caja 14 kilos empedrada'; R = 'Provincia de Limarí'; S = 1071; T = 14 },
    @{ Row = 8; D = 44316; M = 48; N = 14000; O = 14000; P = 14000; Q = '$/caja 14 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 14 },
    @{ Row = 9; D = 44320; M = 45; N = 14000; O = 14000; P = 14000; Q = '$/caja 14 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 14 },
    @{ Row = 10; D = 44314; M = 56; N = 14000; O = 14000; P = 14000; Q = '$/caja 14 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 14 },
    @{ Row = 11; D = 44592; M = 54; N = 20000; O = 20000; P = 20000; Q = '$/caja 15 kilos empedrada'; R = 'Provincia de Limarí'; S = 1333; T = 15 },
    @{ Row = 12; D = 44322; M = 50; N = 14000; O = 14000; P = 14000; Q = '$/caja 14 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 14 },
    @{ Row = 13; D = 44270; M = 85; N = 12000; O = 12000; P = 12000; Q = '$/caja 14 kilos granel'; R = 'Provincia del Elquí'; S = 857; T = 14 },
    @{ Row = 14; D = 44323; M = 60; N = 14000; O = 14000; P = 14000; Q = '$/caja 14 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 14 },
    @{ Row = 15; D = 44260; M = 56; N = 13000; O = 13000; P = 13000; Q = '$/caja 14 kilos empedrada'; R = 'Provincia del Elquí'; S = 929; T = 14 },
    @{ Row = 16; D = 44585; M = 50; N = 22500; O = 22500; P = 22500; Q = '$/caja 15 kilos empedrada'; R = 'Provincia de Limarí'; S = 1500; T = 15 },
    @{ Row = 17; D = 44614; M = 54; N = 14000; O = 14000; P = 14000; Q = '$/caja 14 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 14 },
    @{ Row = 18; D = 44245; M = 50; N = 15000; O = 15000; P = 15000; Q = '$/caja 15 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 15 },
    @{ Row = 19; D = 44627; M = 56; N = 17000; O = 17000; P = 17000; Q = '$/caja 14 kilos empedrada'; R = 'Provincia de Limarí'; S = 1214; T = 14 },
    @{ Row = 20; D = 44252; M = 60; N = 14000; O = 14000; P = 14000; Q = '$/caja 14 kilos empedrada'; R = 'Provincia de Limarí'; S = 1000; T = 14 },
    @{ Row = 21; D = 44271; M = 50; N = 12000; O = 12000; P = 12000; Q = '$/caja 14 kilos granel'; R = 'Provincia del Elquí'; S = 857; T = 14 },
    @{ Row = 22; D = 44238; M = 60; N = 15000; O = 15000; P = 15000; Q = '$/caja 15 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 15 },
    @{ Row = 23; D = 44616; M = 70; N = 14000; O = 14000; P = 14000; Q = '$/caja 14 kilos empedrada'; R = 'Provincia de Limarí'; S = 1000; T = 14 },
    @{ Row = 24; D = 44278; M = 45; N = 13000; O = 13000; P = 13000; Q = '$/caja 14 kilos empedrada'; R = 'Provincia del Elquí'; S = 929; T = 14 },
    @{ Row = 25; D = 44242; M = 45; N = 12000; O = 12000; P = 12000; Q = '$/caja 15 kilos granel'; R = 'Provincia de Limarí'; S = 800; T = 15 },
    @{ Row = 26; D = 44239; M = 70; N = 15000; O = 15000; P = 15000; Q = '$/caja 15 kilos granel'; R = 'Provincia de Limarí'; S = 1000; T = 15 }
)

foreach ($rec in $rows) {
    $r = $rec.Row
    $ws.Cells.Item($r, 4).Value = $rec.D
    $ws.Cells.Item($r, 13).Value = $rec.M
    $ws.Cells.Item($r, 14).Value = $rec.N
    $ws.Cells.Item($r, 15).Value = $rec.O
    $ws.Cells.Item($r, 16).Value = $rec.P
    $ws.Cells.Item($r, 17).Value = $rec.Q
    $ws.Cells.Item($r, 18).Value = $rec.R
    $ws.Cells.Item($r, 19).Value = $rec.S
    $ws.Cells.Item($r, 20).Value = $rec.T
}
